$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Rangahau Hauora Māori"
$ws.Range("C5").Value = "Supporting Māori health research that upholds rangatiratanga and uses and advances Māori knowledge, resources and people"
